$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 0. Grab the cell formats we will need for the new block further down the
#    sheet BEFORE the source cells get cleared out below (PasteSpecial with
#    formats-only re-uses the existing style slot instead of minting a new
#    one, matching how the sheet already looks).
# ---------------------------------------------------------------------------
$ws.Range("I154").Copy()
$ws.Range("I174").PasteSpecial(-4122)
$ws.Range("J154").Copy()
$ws.Range("J174").PasteSpecial(-4122)
$ws.Range("A1").Copy()
$ws.Range("I175").PasteSpecial(-4122)
$ws.Range("A1").Copy()
$ws.Range("J175").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 1. Remove the old "helper calculation" block that lived in columns I/J next
#    to the 2021 week-2 and week-3 data rows (rows 149-155 and 158-160).
# ---------------------------------------------------------------------------
$ws.Range("I149:J154").Clear()
$ws.Range("I155:J155").Clear()
$ws.Range("I158:I160").Clear()

# ---------------------------------------------------------------------------
# 2. Add the new "Source: 2021-02-23" marker row (row 169), styled the same
#    way as the existing "From here from 2021-02-16" marker row (row 99):
#    yellow highlighted cell in column B.
# ---------------------------------------------------------------------------
$ws.Range("B99").Copy()
$ws.Range("B169").PasteSpecial(-4122)
$ws.Range("B169").Value = "Source: 2021-02-23"

# Re-seed the same helper calculation block (now placed after the new data).
$ws.Range("I169").Value = "Länge für Einheit (cm)"
$ws.Range("J169").Formula = "=26.4-6"

# ---------------------------------------------------------------------------
# 3. Add the new week-5 2021 data rows (170-175), mirroring the layout of the
#    existing weekly blocks (A=year, B=week, C=age_group, D=n_tests,
#    E=pct_of_tests_positive).
# ---------------------------------------------------------------------------
$ws.Range("A170").Value = 2021
$ws.Range("B170").Value = 5
$ws.Range("C170").Value = "0-4"
$ws.Range("D170").Value = 7819
$ws.Range("E170").Value = 6.7
$ws.Range("I170").Value = "#:"
$ws.Range("J170").Value = 10

$ws.Range("A171").Value = 2021
$ws.Range("B171").Value = 5
$ws.Range("C171").Value = "5-14"
$ws.Range("D171").Value = 10664
$ws.Range("E171").Value = 9.8
$ws.Range("I171").Value = "# / cm"
$ws.Range("J171").Formula = "=J170 / J169"

$ws.Range("A172").Value = 2021
$ws.Range("B172").Value = 5
$ws.Range("C172").Value = "15-34"
$ws.Range("D172").Value = 95972
$ws.Range("E172").Value = 6.8
$ws.Range("I172").Value = "Achsenabschnitt (cm)"
$ws.Range("J172").Value = 6

$ws.Range("A173").Value = 2021
$ws.Range("B173").Value = 5
$ws.Range("C173").Value = "35-59"
$ws.Range("D173").Value = 153555
$ws.Range("E173").Value = 6.5
$ws.Range("I173").Value = "Achsenabschnitt (#)"
$ws.Range("J173").Value = 5

$ws.Range("A174").Value = 2021
$ws.Range("B174").Value = 5
$ws.Range("C174").Value = "60-79"
$ws.Range("D174").Value = 86730
$ws.Range("E174").Value = 7.3

# I174/J174 already carry the bold-red style copied from I154/J154 above.
$ws.Range("I174").Value = "Gemessene Höhe (cm)"
$ws.Range("J174").Value = 9

$ws.Range("A175").Value = 2021
$ws.Range("B175").Value = 5
$ws.Range("C175").Value = ">=80"
$ws.Range("D175").Value = 53318
$ws.Range("E175").Value = 11.3

# I175/J175 already carry the bold style copied from the header row (A1) above.
$ws.Range("I175").Value = "Zahl"
$ws.Range("J175").Formula = "=(J174-J172)*J171+J173"

# ---------------------------------------------------------------------------
# 4. Re-add the three standalone helper formulas further down the sheet
#    (rows 178-180, leaving rows 176-177 blank exactly like before).
# ---------------------------------------------------------------------------
$ws.Range("I178").Formula = "=50000/55"
$ws.Range("I179").Formula = "=50000 + 909 * 9"
$ws.Range("I180").Formula = "=909*8"

# ---------------------------------------------------------------------------
# 5. Restore the view state (frozen pane selection) to match the new extent
#    of the data.
# ---------------------------------------------------------------------------
$ws.Range("E176").Select()
